$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated result values for rows 2-25 (columns B,C,D,E,G,H,I,K,N,O)
# Columns F, J, L, M remain 0 (unchanged) and column A remains the index (unchanged).
$cols = @("B", "C", "D", "E", "G", "H", "I", "K", "N", "O")

$data = @{}
$data[2] = @(0.3210057619961617, 0.06583081443092453, 0.07913462391887549, 0.4192626335607201, 0.4015654146335876, 0.5554557990681275, 0.4181279530127604, 0.3509234277196356, 1.137998050069106, 1.859820966578127)
$data[3] = @(0.2817881569998804, 0.05844957747930835, 0.07170496381543501, 0.3658471276621782, 0.3998811745049053, 0.5583478614497039, 0.4214918617449577, 0.3062570669490867, 1.151183850461116, 1.862040228453822)
$data[4] = @(0.2576824005849971, 0.05389144221986442, 0.06717585165084472, 0.3331283927742703, 0.3991589699901326, 0.5603905797079562, 0.4238079520904456, 0.2787679261174674, 1.159702404897075, 1.864613814121697)
$data[5] = @(0.2478530525505391, 0.05202745296656985, 0.0653384569556863, 0.3198135427983715, 0.3989429235999182, 0.5612901380011124, 0.4248147362737527, 0.2675502918134214, 1.16327997456812, 1.865966688780901)
$data[6] = @(0.2462205451095372, 0.05171754661849093, 0.06503385831226183, 0.3176036846932391, 0.3989117706286365, 0.5614435637061135, 0.4249857126353263, 0.2656866844114063, 1.163880440006288, 1.866209689312626)
$data[7] = @(0.2575498624480588, 0.05386633013331732, 0.06715103844418024, 0.3329487518788028, 0.3991557396876857, 0.5604024396448608, 0.4238212751059756, 0.2786167036745724, 1.159750223337848, 1.864630828618886)
$data[8] = @(0.3074892005877814, 0.06329118281389867, 0.07656609704184802, 0.4008277555269757, 0.4009198391836151, 0.556397569500561, 0.4192357690242545, 0.3355359739671542, 1.142456822934324, 1.860334597879159)
$data[9] = @(0.4051976049604491, 0.08156622943380398, 0.09528855938839342, 0.5346341580031009, 0.4068631249968746, 0.5506627488847897, 0.4122354427053025, 0.4466332542195914, 1.111897167526845, 1.861537947011016)
$data[10] = @(0.4768334492692361, 0.09486720436534313, 0.1092034622807034, 0.6334747690240192, 0.4127574535367557, 0.5477419707256672, 0.408311175894827, 0.5279265294468587, 1.091489476366649, 1.868322850938313)
$data[11] = @(0.5093870392788062, 0.1008910610760836, 0.1155687003189456, 0.6785797199324151, 0.4157736793962101, 0.5466941654632222, 0.4067916082161993, 0.5648352243764236, 1.082649590442742, 1.872697764035024)
$data[12] = @(0.5217089699976611, 0.1031682681253585, 0.1179841176539895, 0.6956819536012802, 0.416964218192021, 0.5463377945005163, 0.4062544624496525, 0.5788008719985385, 1.079366001163599, 1.8745402102154)
$data[13] = @(0.5190554706768467, 0.102678004780671, 0.1174636907006175, 0.6919976794771259, 0.416705660272811, 0.5464127477664391, 0.4063684424077145, 0.57579361175965, 1.080070338256949, 1.874135136889606)
$data[14] = @(0.5104008833412195, 0.1010784865925132, 0.1157673173435683, 0.6799862817011757, 0.4158706553612461, 0.5466640366066429, 0.4067466491841145, 0.5659844075042031, 1.082378166762725, 1.872845617455823)
$data[15] = @(0.5050989787829394, 0.1000982268784014, 0.1147288950195531, 0.6726318613147555, 0.4153654947018168, 0.5468232215106354, 0.4069832997710954, 0.5599745551974138, 1.083800100433917, 1.87207995569787)
$data[16] = @(0.474705267789119, 0.09447298801630666, 0.1087881860797211, 0.6305300428183358, 0.4125670920681017, 0.5478161009098557, 0.408415835339305, 0.5255129658966382, 1.092076118668656, 1.868062908135869)
$data[17] = @(0.4560506464667355, 0.09101518779843332, 0.1051527658781737, 0.6047393872206044, 0.4109362766112383, 0.5484971530314198, 0.409362739115565, 0.5043530817621331, 1.097266882274919, 1.865928922730859)
$data[18] = @(0.4453178410987562, 0.08902383016274484, 0.1030650909739279, 0.5899185425350879, 0.4100297815494827, 0.5489153118270877, 0.4099323674316473, 0.4921757277098777, 1.100294227809028, 1.864822754279089)
$data[19] = @(0.4416833709607317, 0.08834915866776782, 0.1023588117934366, 0.5849026886857018, 0.4097282629532231, 0.5490614329226418, 0.4101295239011371, 0.4880515416024878, 1.101326404056097, 1.864469033524017)
$data[20] = @(0.4580367937053893, 0.09138353788250697, 0.1055394189225609, 0.6074834630451704, 0.4111066173648794, 0.5484219177290015, 0.409259352133958, 0.5066062897310246, 1.096709994479404, 1.866143537339354)
$data[21] = @(0.5129430964077528, 0.1015484099287391, 0.1162654468073612, 0.6835137116397902, 0.4161146024372187, 0.5465891300906378, 0.4066345210296411, 0.5688659078102489, 1.081698567121339, 1.873219335163952)
$data[22] = @(0.5487957757808886, 0.1081690100076571, 0.1233048862572446, 0.7333326544318766, 0.4196695822985674, 0.5456268489961644, 0.4051422175224708, 0.6094925453009239, 1.072260135428284, 1.878926765900189)
$data[23] = @(0.5296636074293133, 0.1046375640716803, 0.1195451296106, 0.7067310425767346, 0.417746354010859, 0.5461188763742939, 0.405918238125107, 0.5878153363311753, 1.077263500650087, 1.875781347768452)
$data[24] = @(0.4571388814648856, 0.09121701746684607, 0.1053646056772379, 0.6062428461189882, 0.4110295095037344, 0.5484558487100486, 0.409306014769605, 0.5055876525433405, 1.096961629178749, 1.86604613416452)
$data[25] = @(0.3787902538665548, 0.07664445876912396, 0.09019571578306795, 0.4983502696562567, 0.4049879914126961, 0.5519872491312867, 0.4139155030351027, 0.4166355566712241, 1.119805387275221, 1.860178425091277)

foreach ($r in $data.Keys) {
    $rowValues = $data[$r]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $colLetter = $cols[$i]
        $ws.Range("$colLetter$r").Value = $rowValues[$i]
    }
}

Write-Output "Updated $($data.Keys.Count) rows"
